# Corrections following third round of review
#
# The "Materials" sheet (sheet2.xml) had a stray "subgenus" field mapping
# column (header "subgenus", value "${subgenus}") sitting between "genus"
# and "specificEpithet". Remove that entire column so everything to its
# right shifts left by one, and the two now-unused shared strings
# ("subgenus" / "${subgenus}") drop out of the workbook's shared string
# table (and every other shared-string reference that came after them,
# e.g. on the ExternalLinks sheet, is renumbered accordingly).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Materials")
$ws.Columns("AS:AS").Delete()
